$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Brian Hanson (row 9) adds himself to the ADNCS subgroup
$ws.Range("D9").Value = "ADNCS"

# Leave the selection where the editor last clicked, matching the saved view
$ws.Range("D10").Select() | Out-Null
